$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.918.29'
$ws.Range('E2').Value = '  -3.55%  '

$ws.Range('D3').Value = '1.854.98'
$ws.Range('E3').Value = '  -3.03%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.30'
$ws.Range('E5').Value = '  -2.33%  '

$ws.Range('E6').Value = '  -0.06%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4345'
$ws.Range('E7').Value = '  -5.39%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3671'
$ws.Range('E8').Value = '  -3.85%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07476'
$ws.Range('E9').Value = '  -3.17%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9370'
$ws.Range('E10').Value = '  -4.54%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.30'
$ws.Range('E11').Value = '  -3.65%  '

$ws.Range('D12').Value = '1.927.64'
$ws.Range('E12').Value = '  +0.49%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.679'
$ws.Range('E13').Value = '  -3.95%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.415'
$ws.Range('E14').Value = '  -4.51%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06864'
$ws.Range('E15').Value = '  -2.39%  '

$ws.Range('E16').Value = '  +0.00%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '81.36'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009017'
$ws.Range('E18').Value = '  -4.78%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').Value = '  -0.06%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.86'
$ws.Range('E20').Value = '  -4.96%  '

$ws.Range('D21').Value = '27.898.28'
$ws.Range('E21').Value = '  -3.56%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.095'
$ws.Range('E22').Value = '  -4.30%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.03'
$ws.Range('E23').Value = '  +1.31%  '

$ws.Range('D24').Value = '2.134.51'
$ws.Range('E24').Value = '  -0.54%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.007'
$ws.Range('E25').Value = '  -4.11%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.07'
$ws.Range('E26').Value = '  -2.77%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.36'
$ws.Range('E27').Value = '  -3.71%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.367'
$ws.Range('E28').Value = '  -5.34%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '113.19'
$ws.Range('E29').Value = '  -3.63%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.724'
$ws.Range('E30').Value = '  -7.72%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08946'
$ws.Range('E31').Value = '  -3.76%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7997'
$ws.Range('E32').Value = '  -8.19%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.816'
$ws.Range('E33').Value = '  -5.27%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.998'
$ws.Range('E34').Value = '  -5.04%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.166'
$ws.Range('E35').Value = '  -7.06%  '

$ws.Range('E36').Value = '  -0.05%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.115'
$ws.Range('E37').Value = '  -4.06%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05412'
$ws.Range('E38').Value = '  -5.56%  '

$ws.Range('E39').Value = '  -3.99%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.919'
$ws.Range('E40').Value = '  +2.13%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5224'
$ws.Range('E41').Value = '  -4.96%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.968'
$ws.Range('E42').Value = '  -6.10%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1677'
$ws.Range('E43').Value = '  -4.56%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.720'
$ws.Range('E44').Value = '  -6.60%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.06716'
$ws.Range('E45').Value = '  -2.80%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4858'
$ws.Range('E46').Value = '  -6.30%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.62'
$ws.Range('E47').Value = '  -5.39%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '106.29'
$ws.Range('E48').Value = '  -3.96%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.920'
$ws.Range('E49').Value = '  -8.60%  '

$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.000'
$ws.Range('E50').Value = '  -0.06%  '

$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.671'
$ws.Range('E51').Value = '  -6.23%  '
